# ControlStateDiagram.xlsx - "cu and control state diagram update"
#
# The author unhid the (previously hidden) control-unit bit columns B:E on
# Sheet1 so the CU/control-state bit values are visible again, and moved
# the viewport/selection down to the part of the state diagram they were
# working on (around row ~111, cell C106) instead of the old position
# (around row ~94, cell P149).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Unhide the control-unit columns B, C, D, E (were hidden="1" in the xml).
$ws.Columns("B:E").Hidden = $false

# Move the selection to where the author was working, and scroll the
# viewport so that area is visible (row 111 at the top of the window).
$ws.Range("C106").Select()
$excel.ActiveWindow.ScrollRow = 111
$excel.ActiveWindow.ScrollColumn = 1
